$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cells already carrying the fill styles we need to reuse
# (copying preserves the exact existing style index instead of minting a
# new one, which matters for the gray/theme-based "style 1" fill).
$styleSrc1 = $ws.Range("C3")   # plain grey fill (s="1")
$styleSrc2 = $ws.Range("E5")   # red fill (s="2")
$styleSrc3 = $ws.Range("H2")   # purple fill (s="3")
$styleSrc5 = $ws.Range("F5")   # green fill (s="5")

function Paint($rangeAddr, $src) {
    $src.Copy()
    $ws.Range($rangeAddr).PasteSpecial(-4122)
}

# Row 13 - A=4
$ws.Range("A13").Value = 4
Paint "B13" $styleSrc1
Paint "C13:F13" $styleSrc3
Paint "G13" $styleSrc5
Paint "H13" $styleSrc2
Paint "I13:V13" $styleSrc1

# Row 14 - A=3
$ws.Range("A14").Value = 3
Paint "B14:E14" $styleSrc1
Paint "F14" $styleSrc2
Paint "G14:V14" $styleSrc1

# Row 15 - A=2
$ws.Range("A15").Value = 2
Paint "B15:D15" $styleSrc1
Paint "E15" $styleSrc2
Paint "F15:V15" $styleSrc1

# Row 16 - A=1
$ws.Range("A16").Value = 1
Paint "B16:C16" $styleSrc1
Paint "D16" $styleSrc2
Paint "E16:V16" $styleSrc1

# Row 17 - A=0
$ws.Range("A17").Value = 0
Paint "B17:C17" $styleSrc2
Paint "D17:V17" $styleSrc1

$excel.CutCopyMode = 0
$ws.Range("D17").Select()
